$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "C1qtnf1"
$ws.Range("C2").Value = "Avpr2"
$ws.Range("D2").Value = "Inflammatory-Mac"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.937483333333333
$ws.Range("H2").Value = 14.81245
$ws.Range("I2").Value = 0.2144582670332307
$ws.Range("J2").Value = 0.2144582670332307
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.05743166666666667
$ws.Range("N2").Value = 0.172295
$ws.Range("O2").Value = 0.4025715794441874
$ws.Range("P2").Value = 0.4025715794441875
$ws.Range("Q2").Value = 0.2835678969722222
$ws.Range("R2").Value = 2.55211107275
$ws.Range("S2").Value = 0.08633480328443098
$ws.Range("T2").Value = 0.08633480328443101

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "C1qtnf1"
$ws.Range("C3").Value = "Avpr2"
$ws.Range("D3").Value = "Resolving-Mac"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.937483333333333
$ws.Range("H3").Value = 14.81245
$ws.Range("I3").Value = 0.2144582670332307
$ws.Range("J3").Value = 0.2144582670332307
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.08523033333333334
$ws.Range("N3").Value = 0.255691
$ws.Range("O3").Value = 0.5974284205558126
$ws.Range("P3").Value = 0.5974284205558126
$ws.Range("Q3").Value = 0.4208233503277777
$ws.Range("R3").Value = 3.78741015295
$ws.Range("S3").Value = 0.1281234637487997
$ws.Range("T3").Value = 0.1281234637487997

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "C1qtnf1"
$ws.Range("C4").Value = "Avpr2"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 12.04564066666667
$ws.Range("H4").Value = 36.136922
$ws.Range("I4").Value = 0.5231991782611944
$ws.Range("J4").Value = 0.5231991782611944
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.05743166666666667
$ws.Range("N4").Value = 0.172295
$ws.Range("O4").Value = 0.4025715794441874
$ws.Range("P4").Value = 0.4025715794441875
$ws.Range("Q4").Value = 0.6918012195544444
$ws.Range("R4").Value = 6.22621097599
$ws.Range("S4").Value = 0.21062511955651
$ws.Range("T4").Value = 0.21062511955651

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "C1qtnf1"
$ws.Range("C5").Value = "Avpr2"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 12.04564066666667
$ws.Range("H5").Value = 36.136922
$ws.Range("I5").Value = 0.5231991782611944
$ws.Range("J5").Value = 0.5231991782611944
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.08523033333333334
$ws.Range("N5").Value = 0.255691
$ws.Range("O5").Value = 0.5974284205558126
$ws.Range("P5").Value = 0.5974284205558126
$ws.Range("Q5").Value = 1.026653969233555
$ws.Range("R5").Value = 9.239885723101999
$ws.Range("S5").Value = 0.3125740587046844
$ws.Range("T5").Value = 0.3125740587046844

# Row 6
$ws.Range("A6").Value = "Inflammatory-Mac"
$ws.Range("B6").Value = "C1qtnf1"
$ws.Range("C6").Value = "Avpr2"
$ws.Range("D6").Value = "Inflammatory-Mac"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.065339666666667
$ws.Range("H6").Value = 6.196019
$ws.Range("I6").Value = 0.08970747561983136
$ws.Range("J6").Value = 0.08970747561983136
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.05743166666666667
$ws.Range("N6").Value = 0.172295
$ws.Range("O6").Value = 0.4025715794441874
$ws.Range("P6").Value = 0.4025715794441875
$ws.Range("Q6").Value = 0.1186158992894444
$ws.Range("R6").Value = 1.067543093605
$ws.Range("S6").Value = 0.03611368014822644
$ws.Range("T6").Value = 0.03611368014822645

# Row 7
$ws.Range("A7").Value = "Inflammatory-Mac"
$ws.Range("B7").Value = "C1qtnf1"
$ws.Range("C7").Value = "Avpr2"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.065339666666667
$ws.Range("H7").Value = 6.196019
$ws.Range("I7").Value = 0.08970747561983136
$ws.Range("J7").Value = 0.08970747561983136
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.08523033333333334
$ws.Range("N7").Value = 0.255691
$ws.Range("O7").Value = 0.5974284205558126
$ws.Range("P7").Value = 0.5974284205558126
$ws.Range("Q7").Value = 0.1760295882365556
$ws.Range("R7").Value = 1.584266294129
$ws.Range("S7").Value = 0.05359379547160491
$ws.Range("T7").Value = 0.05359379547160491

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "C1qtnf1"
$ws.Range("C8").Value = "Avpr2"
$ws.Range("D8").Value = "Inflammatory-Mac"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.956235666666667
$ws.Range("H8").Value = 5.868707000000001
$ws.Range("I8").Value = 0.0849685725822393
$ws.Range("J8").Value = 0.0849685725822393
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.05743166666666667
$ws.Range("N8").Value = 0.172295
$ws.Range("O8").Value = 0.4025715794441874
$ws.Range("P8").Value = 0.4025715794441875
$ws.Range("Q8").Value = 0.1123498747294445
$ws.Range("R8").Value = 1.011148872565
$ws.Range("S8").Value = 0.03420593246755015
$ws.Range("T8").Value = 0.03420593246755016

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "C1qtnf1"
$ws.Range("C9").Value = "Avpr2"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.956235666666667
$ws.Range("H9").Value = 5.868707000000001
$ws.Range("I9").Value = 0.0849685725822393
$ws.Range("J9").Value = 0.0849685725822393
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.08523033333333334
$ws.Range("N9").Value = 0.255691
$ws.Range("O9").Value = 0.5974284205558126
$ws.Range("P9").Value = 0.5974284205558126
$ws.Range("Q9").Value = 0.1667306179485556
$ws.Range("R9").Value = 1.500575561537
$ws.Range("S9").Value = 0.05076264011468915
$ws.Range("T9").Value = 0.05076264011468915

# Row 10
$ws.Range("A10").Value = "Resolving-Mac"
$ws.Range("B10").Value = "C1qtnf1"
$ws.Range("C10").Value = "Avpr2"
$ws.Range("D10").Value = "Inflammatory-Mac"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 2.018350333333334
$ws.Range("H10").Value = 6.055051000000001
$ws.Range("I10").Value = 0.08766650650350422
$ws.Range("J10").Value = 0.0876665065035042
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.05743166666666667
$ws.Range("N10").Value = 0.172295
$ws.Range("O10").Value = 0.4025715794441874
$ws.Range("P10").Value = 0.4025715794441875
$ws.Range("Q10").Value = 0.1159172235605556
$ws.Range("R10").Value = 1.043255012045
$ws.Range("S10").Value = 0.03529204398746982
$ws.Range("T10").Value = 0.03529204398746982

# Row 11
$ws.Range("A11").Value = "Resolving-Mac"
$ws.Range("B11").Value = "C1qtnf1"
$ws.Range("C11").Value = "Avpr2"
$ws.Range("D11").Value = "Resolving-Mac"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 2.018350333333334
$ws.Range("H11").Value = 6.055051000000001
$ws.Range("I11").Value = 0.08766650650350422
$ws.Range("J11").Value = 0.0876665065035042
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.08523033333333334
$ws.Range("N11").Value = 0.255691
$ws.Range("O11").Value = 0.5974284205558126
$ws.Range("P11").Value = 0.5974284205558126
$ws.Range("Q11").Value = 0.1720246716934445
$ws.Range("R11").Value = 1.548222045241
$ws.Range("S11").Value = 0.05237446251603439
$ws.Range("T11").Value = 0.05237446251603439
